$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel auto-detection; force text
# so the stored cell stays a string (matching the source data which is
# pre-formatted price text), then restore default styling.
$textCells = @("D5","D6","D7","D8","D11","D12","D13","D14","D16","D20","D21","D22","D23","D24","D25","D28","D29","D31","D32","D33","D37","D38","D40","D41","D43","D45","D46","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.719.32'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '3.317.84'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '581.64'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '173.93'
$ws.Range("E6").Value = '  -7.18%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.584'
$ws.Range("E8").Value = '  -1.95%  '
$ws.Range("D9").Value = '3.314.13'
$ws.Range("E10").Value = '  -5.35%  '
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").Value = '45.28'
$ws.Range("E12").Value = '  -4.94%  '
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").Value = '663.10'
$ws.Range("E14").Value = '  +3.53%  '
$ws.Range("D15").Value = '3.858.47'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '8.37'
$ws.Range("E16").Value = '  -3.20%  '
$ws.Range("D17").Value = '67.849.34'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D19").Value = '3.318.16'
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").Value = '17.41'
$ws.Range("E20").Value = '  -3.80%  '
$ws.Range("D21").Value = '10.87'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").Value = '0.886'
$ws.Range("E22").Value = '  -2.87%  '
$ws.Range("D23").Value = '5.39'
$ws.Range("E23").Value = '  +5.34%  '
$ws.Range("D24").Value = '16.86'
$ws.Range("E24").Value = '  -6.38%  '
$ws.Range("D25").Value = '97.42'
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("E26").Value = '  -5.31%  '
$ws.Range("E27").Value = '  -6.94%  '
$ws.Range("D28").Value = '9.23'
$ws.Range("E28").Value = '  -6.03%  '
$ws.Range("D29").Value = '33.41'
$ws.Range("E29").Value = '  +2.09%  '
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("D32").Value = '582.99'
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("D33").Value = '10.93'
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").Value = '3.715.48'
$ws.Range("E35").Value = '  -7.78%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = '56.81'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").Value = '3.29'
$ws.Range("E38").Value = '  -14.93%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = '32.35'
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  -7.49%  '
$ws.Range("E42").Value = '  -5.48%  '
$ws.Range("D43").Value = '0.331'
$ws.Range("E43").Value = '  -3.82%  '
$ws.Range("D44").Value = '0.0₃0661'
$ws.Range("E44").Value = '  -6.34%  '
$ws.Range("D45").Value = '3.26'
$ws.Range("E45").Value = '  -4.47%  '
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  -4.45%  '
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  -2.88%  '
$ws.Range("D51").Value = '127.35'
$ws.Range("E51").Value = '  -0.51%  '

# Clear the explicit format again so the cell style matches the original
# (no explicit style index), now that the value is safely stored as text.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
